$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.897.07"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "2.650.62"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.14%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  -0.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.57"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.03%  "
$ws.Range("E10").Value = "  +0.89%  "
$ws.Range("E11").Value = "  +3.79%  "
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").Value = "3.119.44"
$ws.Range("E13").Value = "  +1.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +11.98%  "
$ws.Range("D15").Value = "60.867.66"
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("D17").Value = "2.660.47"
$ws.Range("E17").Value = "  +1.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "350.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.532"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.27%  "
$ws.Range("E25").Value = "  +0.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.82%  "
$ws.Range("E28").Value = "  +8.38%  "
$ws.Range("D29").Value = "0.0₃0812"
$ws.Range("E29").Value = "  +2.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.72%  "
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "166.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.48%  "
$ws.Range("E33").Value = "  +1.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.56"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.81%  "
$ws.Range("E35").Value = "  +5.26%  "
$ws.Range("E36").Value = "  +8.10%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "341.14"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +14.41%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.66"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.904"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.59"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.40"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "133.97"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.81%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.619"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0562"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.53%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0249"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.68%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.74%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0995"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.02%  "
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("D51").Value = "2.107.89"
$ws.Range("E51").Value = "  +4.15%  "
